# Rename the three embedded logo pictures (two Pearson logos in the
# footers, one BTEC logo in the header) as captured by the diff:
#   footer1 Pearson logo (id=3):  image2.png -> image1.png
#   footer2 Pearson logo (id=2):  image2.png -> image1.png
#   header1 BTEC logo   (id=1):  image1.jpg -> image2.jpg
#
# The pictures themselves, their ids, sizes and alt-text/description are
# left untouched - only the shape's Name changes.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Headers ------------------------------------------------------
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $shp = $hdr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers --------------------------------------------------------
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $shp = $ftr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
